# Execute Tests on Android mobile browser
# Adds a "Locator Type" column (D) to Sheet1 classifying each existing
# locator value (column C) as an "Xpath" or "CSS" locator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for the new column
$ws.Range("D1").Value = "Locator Type"

# Locator-type values for rows 2..14, derived from the existing C column
# (values starting with "//" are XPath expressions, everything else is CSS).
$locatorTypes = @(
    "Xpath",  # row 2  - //div[@aria-label= 'Cookie banner']...
    "CSS",    # row 3  - input#name
    "CSS",    # row 4  - #name + span.error-msg
    "CSS",    # row 5  - input#email
    "CSS",    # row 6  - #email + span.error-msg
    "CSS",    # row 7  - input#password
    "CSS",    # row 8  - input#orgName
    "CSS",    # row 9  - #orgName + span.error-msg
    "Xpath",  # row 10 - //input[@type='checkbox' and @name='termsOfServiceAccepted']
    "Xpath",  # row 11 - //input[@type='checkbox' and @name='marketingOptIn']
    "CSS",    # row 12 - button#formSubmit
    "Xpath",  # row 13 - //header[contains(@class,'site-header')]...
    "CSS"     # row 14 - img[src*='successful-res-icon.svg']
)

for ($i = 0; $i -lt $locatorTypes.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $locatorTypes[$i]
}

# Widen column C (now holding longer locator values) and size the new
# column D to fit its "Locator Type"/"Xpath"/"CSS" contents.
$ws.Columns.Item(3).ColumnWidth = 59.6
$ws.Columns.Item(4).ColumnWidth = 9.7

# Move the active selection, matching the saved state after editing.
$ws.Range("C14").Select()
